$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXmlByAnchor($anchorText, $innerXml) {
    foreach ($par in $d.Paragraphs) {
        $rngText = $par.Range.Text
        if ($rngText.StartsWith($anchorText)) {
            $xml = "<w:p $wns>$innerXml</w:p>"
            $par.Range.InsertXML($xml)
            return
        }
    }
    throw ("Paragraph starting with " + $anchorText + " not found")
}

# --- pt_programa ---
Set-ParagraphXmlByAnchor 'Princípios elementares em química: Sistemas' '<w:r><w:t>Princípios elementares em química: Sistemas de Unidades (Definição das Unidades mais usadas em Engenharia e transformações entre sistemas).</w:t><w:br/><w:t>Estrutura Atômica e Tabela Periódica: Natureza elétrica da matéria. A carga do elétron. O núcleo do átomo. Espectros de emissão e de absorção atômica. Configuração eletrônica dos elementos. Partículas Elementares. A Lei e a tabela Periódica.</w:t><w:br/><w:t xml:space="preserve">A Ligação Química: A ligação eletrovalente. A ligação covalente. Hibridação. Polaridade da ligação. </w:t><w:br/><w:t>Natureza dos Compostos: Ácidos e bases (Arrhenius, Bronsted-Lowry e Lewis). Forças intermoleculares.</w:t><w:br/><w:t>Reações Químicas em Solução Aquosa : Terminologia das soluções. Eletrólitos e não eletrólitos. Reações iônicas. Reações sem transferência de elétron e seu balanceamento. Preparação de sais inorgânicos (por dupla troca). Oxidação e redução. Número de oxidação. Reações de óxido redução. Métodos de balanceamento de reações de oxi-redução (Variação do Nox, via decomposição do agente oxidante, íon-elétron e pelo Potencial Padrão de Redução).</w:t><w:br/><w:t>Gases: Variáveis de estado. Lei combinada dos gases. Experiência de Torriceli. Teoria cinética dos gases. Gás ideal e real. Princípio de Avogadro.</w:t><w:br/><w:t>Soluções: Natureza das soluções. Dispersões coloidais e suspensões. Tipos de soluções. Unidades de concentração (Molaridade, fração molar, ppm, normalidade, molalidade). O processo de dissolução. Calor de dissolução. Solubilidade e temperatura.</w:t><w:br/><w:t>Estequiometria e Cálculos em Química : Cálculos baseados em equações químicas. Cálculos com reagentes limitantes e reagentes com pureza. Rendimento teórico e centesimal. Resolução de exercícios envolvendo estequiometria industrial.</w:t></w:r>'

# --- en_programa ---
Set-ParagraphXmlByAnchor 'Elementary principles of chemistry: Units' '<w:r><w:rPr><w:i/></w:rPr><w:t>Elementary principles of chemistry: Units Systems (Definition of the most used units in Engineering and transformations between systems).</w:t><w:br/><w:t>Atomic structure and the Periodic Table: electrical nature of matter. The electron charge. The nucleus of the atom. Emission spectra and atomic absorption. Electronic configuration of the elements. Elementary Particles. The Law and the Periodic Table.</w:t><w:br/><w:t xml:space="preserve">The Chemical Bonding: The ionic bonding. The covalent bond. Hybridization. Polarity of covalent bonding. </w:t><w:br/><w:t>Nature of the Compounds: Acids and bases (Arrhenius, Bronsted-Lowry and Lewis). Intermolecular forces.</w:t><w:br/><w:t>Chemical Reactions in Aqueous Solution: Terminology in Solutions. Electrolytes and non electrolytes. Ionic reactions. Reactions without electron transfer and its balancing. Preparation of inorganic salts (metathesis). Oxidation and reduction. Oxidation number. Redox reactions. Redox reactions balancing methods (Variation of Nox,  decomposition of the oxidizing agent, ion-electron and using the Standard Potential of Reduction).</w:t><w:br/><w:t>Gases: State variables. Combined gas law. Experience Torriceli. Kinetic theory of gases. Ideal and real gas. Avogadro''s Principle.</w:t><w:br/><w:t>Solutions: Nature of solutions. Colloidal dispersions and suspensions. Types of solutions. Concentration units (Molarity, mole fraction, ppm, normality, molality). The dissolution process. Heat dissolution. Solubility and temperature.</w:t><w:br/><w:t>Stoichiometric calculations in Chemistry: Calculations based on chemical equations. Calculations with limiting reagents and reagent purity. Theoretical and centesimal yields. Solving of exercises with industrial stoichiometric approach.</w:t></w:r>'

# --- biblio ---
Set-ParagraphXmlByAnchor 'BROWN, T.L. ET al.' '<w:r><w:t>BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007</w:t><w:br/><w:t>ATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006</w:t><w:br/><w:t>BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981</w:t><w:br/><w:t>CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. AMGH Editora Ltda., 2010.</w:t><w:br/><w:t>RUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill</w:t></w:r>'
